# Commit: "test P7 with -10 percent"
# Applies value updates + row deletions across the result sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "general": update scalar summary stats
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value  = 198.2574223937857
$ws.Range("B4").Value  = 0.01199984550476074
$ws.Range("B6").Value  = 30.16742239378572
$ws.Range("B7").Value  = 0
$ws.Range("B8").Value  = 0
$ws.Range("B9").Value  = 0
$ws.Range("B10").Value = 168.09

# ---------------------------------------------------------------
# Sheet "x": reassigned j values for several rows
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value  = 2
$ws.Range("B4").Value  = 6
$ws.Range("B6").Value  = 1
$ws.Range("B7").Value  = 13
$ws.Range("B8").Value  = 11
$ws.Range("B10").Value = 3
$ws.Range("B13").Value = 8
$ws.Range("B14").Value = 12

# ---------------------------------------------------------------
# Sheet "U": reassigned t values for several rows
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("U")
$ws.Range("B2").Value  = 3
$ws.Range("B3").Value  = 2
$ws.Range("B6").Value  = 2
$ws.Range("B8").Value  = 3
$ws.Range("B10").Value = 3
$ws.Range("B12").Value = 3

# ---------------------------------------------------------------
# Sheet "TBar": updated TBar series
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value  = 20
$ws.Range("B4").Value  = 10
$ws.Range("B5").Value  = 20.34885527085025
$ws.Range("B6").Value  = 24.76592070603971
$ws.Range("B7").Value  = 10
$ws.Range("B8").Value  = 20
$ws.Range("B9").Value  = 20.60033324079215
$ws.Range("B10").Value = 22.45367071955468
$ws.Range("B11").Value = 20
$ws.Range("B12").Value = 22.61192465059683
$ws.Range("B13").Value = 30
$ws.Range("B14").Value = 26.71671453559703
$ws.Range("B15").Value = 22.66758337047728

# ---------------------------------------------------------------
# Sheet "y": all data rows removed (header-only sheet now)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("y")
$ws.Range("A2:D11").EntireRow.Delete()

# ---------------------------------------------------------------
# Sheet "Q": updated Q values for rows 7-71
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value  = 177.7200000000015
$ws.Range("C8").Value  = 191.3600000000015
$ws.Range("C9").Value  = 187.3050000000015
$ws.Range("C10").Value = 192.2
$ws.Range("C11").Value = 188.4700000000015
$ws.Range("C12").Value = 67.77500000000072
$ws.Range("C13").Value = 73.77000000000072
$ws.Range("C14").Value = 74.03500000000074
$ws.Range("C15").Value = 73.04500000000073
$ws.Range("C16").Value = 72.66500000000073
$ws.Range("C17").Value = 154.3
$ws.Range("C18").Value = 148.3449999999993
$ws.Range("C19").Value = 128.7049999999993
$ws.Range("C20").Value = 146.3249999999993
$ws.Range("C21").Value = 134.2149999999993
$ws.Range("C22").Value = 288.6699999999997
$ws.Range("C23").Value = 318.1
$ws.Range("C24").Value = 289.5099999999996
$ws.Range("C25").Value = 305.3
$ws.Range("C26").Value = 289.8849999999996
$ws.Range("C27").Value = 62.63000000000022
$ws.Range("C28").Value = 70.92000000000021
$ws.Range("C29").Value = 67.65500000000021
$ws.Range("C30").Value = 66.84500000000021
$ws.Range("C31").Value = 65.41000000000022
$ws.Range("C32").Value = 107.3799999999999
$ws.Range("C33").Value = 112.2399999999999
$ws.Range("C34").Value = 93.78999999999985
$ws.Range("C35").Value = 108.8349999999998
$ws.Range("C36").Value = 94.77999999999986
$ws.Range("C37").Value = 260.7950000000024
$ws.Range("C38").Value = 276.0800000000024
$ws.Range("C39").Value = 265.2800000000024
$ws.Range("C40").Value = 281.9700000000025
$ws.Range("C41").Value = 270.3250000000024
$ws.Range("C42").Value = 216.0400000000002
$ws.Range("C43").Value = 238.0450000000001
$ws.Range("C44").Value = 207.1900000000002
$ws.Range("C45").Value = 221.3450000000001
$ws.Range("C46").Value = 209.0700000000002
$ws.Range("C47").Value = 148.1150000000007
$ws.Range("C48").Value = 158.8850000000007
$ws.Range("C49").Value = 147.3700000000007
$ws.Range("C50").Value = 157.5750000000007
$ws.Range("C51").Value = 150.1400000000007
$ws.Range("C52").Value = 318.7450000000017
$ws.Range("C53").Value = 334.7600000000017
$ws.Range("C54").Value = 327.0100000000016
$ws.Range("C55").Value = 342.6250000000018
$ws.Range("C56").Value = 323.2400000000017
$ws.Range("C57").Value = 288.6699999999997
$ws.Range("C58").Value = 318.1
$ws.Range("C59").Value = 289.5099999999996
$ws.Range("C60").Value = 305.3
$ws.Range("C61").Value = 289.8849999999996
$ws.Range("C62").Value = 318.7450000000017
$ws.Range("C63").Value = 334.7600000000017
$ws.Range("C64").Value = 327.0100000000016
$ws.Range("C65").Value = 342.6250000000018
$ws.Range("C66").Value = 323.2400000000017
$ws.Range("C67").Value = 154.3
$ws.Range("C68").Value = 148.3449999999993
$ws.Range("C69").Value = 128.7049999999993
$ws.Range("C70").Value = 146.3249999999993
$ws.Range("C71").Value = 134.2149999999993

# ---------------------------------------------------------------
# Sheet "R": rows 7-11 zeroed out
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R")
$ws.Range("C7").Value  = 0
$ws.Range("C8").Value  = 0
$ws.Range("C9").Value  = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0

# ---------------------------------------------------------------
# Sheet "L": values shifted between blocks
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("L")
$ws.Range("C2").Value  = 0
$ws.Range("C3").Value  = 0
$ws.Range("C4").Value  = 0
$ws.Range("C5").Value  = 0
$ws.Range("C6").Value  = 0
$ws.Range("C7").Value  = 13.7
$ws.Range("C8").Value  = 6.91
$ws.Range("C9").Value  = 10.68
$ws.Range("C10").Value = 7.39
$ws.Range("C11").Value = 14.68
$ws.Range("C22").Value = 7.25
$ws.Range("C23").Value = 5.4
$ws.Range("C24").Value = 4.755
$ws.Range("C25").Value = 5.8
$ws.Range("C26").Value = 7.48
$ws.Range("C32").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("C34").Value = 0
$ws.Range("C35").Value = 0
$ws.Range("C36").Value = 0
$ws.Range("C42").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("C44").Value = 0
$ws.Range("C45").Value = 0
$ws.Range("C46").Value = 0

# ---------------------------------------------------------------
# Sheet "rho": all data rows removed (header-only sheet now)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2:C11").EntireRow.Delete()

# ---------------------------------------------------------------
# Sheet "alpha": all data rows removed (header-only sheet now)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("alpha")
$ws.Range("A2:C6").EntireRow.Delete()
